$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.012.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.56%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.597.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.36%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("E5").Value = "  -0.46%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.68%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3775"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.52%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3639"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.89"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.250"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.98%  "

# Row 11
$ws.Range("E11").Value = "  -0.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08122"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.574"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.343"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.89%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001241"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.599.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06804"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.65%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.63%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.503"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.83%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "

# Row 23
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "

# Row 24
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.018.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.370"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.790"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.49%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.76%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.45%  "

# Row 29
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.226"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.60%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "

# Row 31
$ws.Range("B31").Value = "WEMIXTOKEN"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.369"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.42%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.787"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -14.23%  "

# Row 33
$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.774.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.41%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9547"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07560"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.22%  "

# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.186"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.58%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.39%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2515"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.54%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.84%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.364"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.18%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.67%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.99%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.22%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6596"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.50%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.991"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.274"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.89%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.52%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07887"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "

# Row 51
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.210"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.14%  "
